$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 "Distance Measurements": insert a new "Measured1" column ---
# before: A..F = Joint1,Joint2,Joint3,Measured,matlab
# after:  A..G = Joint1,Joint2,Joint3,Measured1,Measured2,matlab
$ws1.Columns.Item(5).Insert()

$ws1.Range("E1").Value = "Measured2"
$ws1.Range("F1").Value = "Measured1"

$measured1 = @(0.212, 0.115, 0.233, 0.118, 0.087, 0.136, 0.264, 0.175, 0.231)
for ($i = 0; $i -lt $measured1.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 5).Value = $measured1[$i]
}

# New column + the old "Measured" column (now F) both get the highlighted
# (yellow fill + border) style.
$fmtRange = $ws1.Range("E2:F10")
$fmtRange.Interior.Color = 65535
$fmtRange.Borders.LineStyle = 1

# Column widths for the two new columns.
$ws1.Range("E1:F1").ColumnWidth = 9.83

# --- Sheet2 "DH parameters meaurement": fill in new measurement row ---
$ws2.Range("E4").Value = 0

# --- View / selection state ---
$ws2.Range("I8").Select()

$ws1.Range("E2:E10").Select()
$ws1.Activate()
